$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 10 to make room for two additional children records
$ws.Range("A10:A11").EntireRow.Insert()

# Row 4
$ws.Range("A4:B4").NumberFormat = "@"
$ws.Range("A4").Value = "nChildren"
$ws.Range("B4").Value = "6"
$ws.Range("A4:B4").ClearFormats()

# Row 5
$ws.Range("A5:B5").NumberFormat = "@"
$ws.Range("A5").Value = "timeOfStart"
$ws.Range("B5").Value = "7:00:00"
$ws.Range("A5:B5").ClearFormats()

# Row 6
$ws.Range("A6:H6").NumberFormat = "@"
$ws.Range("A6").Value = "0"
$ws.Range("B6").Value = "16"
$ws.Range("C6").Value = "Collette  "
$ws.Range("D6").Value = "Billi  "
$ws.Range("E6").Value = "-6,9"
$ws.Range("F6").Value = "Elias(mother): 0578741979"
$ws.Range("G6").Value = "7:00:00"
$ws.Range("H6").Value = "20.0"
$ws.Range("A6:H6").ClearFormats()

# Row 7
$ws.Range("A7:H7").NumberFormat = "@"
$ws.Range("A7").Value = "1"
$ws.Range("B7").Value = "14"
$ws.Range("C7").Value = "Lorinda  "
$ws.Range("D7").Value = "Tyron  "
$ws.Range("E7").Value = "-7,8"
$ws.Range("F7").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G7").Value = "7:02:00"
$ws.Range("H7").Value = "18.0"
$ws.Range("A7:H7").ClearFormats()

# Row 8
$ws.Range("A8:H8").NumberFormat = "@"
$ws.Range("A8").Value = "2"
$ws.Range("B8").Value = "14"
$ws.Range("C8").Value = "Lorinda  "
$ws.Range("D8").Value = "Tyron  "
$ws.Range("E8").Value = "-7,8"
$ws.Range("F8").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G8").Value = "7:02:00"
$ws.Range("H8").Value = "18.0"
$ws.Range("A8:H8").ClearFormats()

# Row 9
$ws.Range("A9:H9").NumberFormat = "@"
$ws.Range("A9").Value = "3"
$ws.Range("B9").Value = "19"
$ws.Range("C9").Value = "Jeanine  "
$ws.Range("D9").Value = "Janee  "
$ws.Range("E9").Value = "-7,7"
$ws.Range("F9").Value = "Teresa(mother): 0517627420"
$ws.Range("G9").Value = "7:04:00"
$ws.Range("H9").Value = "16.0"
$ws.Range("A9:H9").ClearFormats()

# Row 10
$ws.Range("A10:H10").NumberFormat = "@"
$ws.Range("A10").Value = "4"
$ws.Range("B10").Value = "20"
$ws.Range("C10").Value = "Ron"
$ws.Range("D10").Value = "Cohen"
$ws.Range("E10").Value = "-8,7"
$ws.Range("F10").Value = "Bernardine(mother): 0576270618"
$ws.Range("G10").Value = "7:06:00"
$ws.Range("H10").Value = "14.0"
$ws.Range("A10:H10").ClearFormats()

# Row 11
$ws.Range("A11:H11").NumberFormat = "@"
$ws.Range("A11").Value = "5"
$ws.Range("B11").Value = "15"
$ws.Range("C11").Value = "Nubia  "
$ws.Range("D11").Value = "Royce  "
$ws.Range("E11").Value = "-9,7"
$ws.Range("F11").Value = "Augustus(father): 0517389040"
$ws.Range("G11").Value = "7:08:00"
$ws.Range("H11").Value = "12.0"
$ws.Range("A11:H11").ClearFormats()

# Row 12
$ws.Range("A12:G12").NumberFormat = "@"
$ws.Range("A12").Value = "school"
$ws.Range("B12").Value = "3"
$ws.Range("C12").Value = "Ironiah"
$ws.Range("D12").Value = "mySchool"
$ws.Range("E12").Value = "0,0"
$ws.Range("F12").Value = "Shir(secretary): 0523345098"
$ws.Range("G12").Value = "7:20:00"
$ws.Range("A12:G12").ClearFormats()

# Row 13
$ws.Range("A13:B13").NumberFormat = "@"
$ws.Range("A13").Value = "cost"
$ws.Range("B13").Value = "65.0"
$ws.Range("A13:B13").ClearFormats()

# Row 14
$ws.Range("A14:B14").NumberFormat = "@"
$ws.Range("A14").Value = "time"
$ws.Range("B14").Value = "20.0"
$ws.Range("A14:B14").ClearFormats()
